$d = $word.ActiveDocument

# 1) Replace the two representative names.
$d.Content.Find.Execute("3.1-Uno es Javier Jiménez Representante del Sistema FSC.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.1-Uno es Betzabet Marín Representante del Sistema FSC.", 2)

$d.Content.Find.Execute("3.2-Otro es Fernando Gómez Representante de Higiene y Seguridad.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.2-Otro es Araceli Becerril Representante de Higiene y Seguridad.", 2)

# 2) Merge " en" + "cuentra al final de la cadena " into one run, removing the old bookmark location.
$d.Content.Find.Execute(" encuentra al final de la cadena ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0) | Out-Null
